# Scheduled runner update: refresh market-price derived columns (H-N)
# across the ALC/ARM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# Sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 2432.0833
$ws.Range("I41").Value = 3949.5
$ws.Range("J41").Value = 2128.6
$ws.Range("K41").Value = 3949.5
$ws.Range("L41").Value = 2128.6
$ws.Range("M41").Value = -3509.5
$ws.Range("N41").Value = -3008.6

# Sheet ALC, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 5907.8184
$ws.Range("I53").Value = 7661
$ws.Range("J53").Value = 1232.6666
$ws.Range("K53").Value = 7661
$ws.Range("L53").Value = 1232.6666
$ws.Range("M53").Value = -7024
$ws.Range("N53").Value = -2506.6666

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 115310.22
$ws.Range("I62").Value = 253549.25
$ws.Range("J62").Value = 4719
$ws.Range("K62").Value = 253549.25
$ws.Range("L62").Value = 4719
$ws.Range("M62").Value = -252925.25
$ws.Range("N62").Value = -5967

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 115310.22
$ws.Range("I65").Value = 253549.25
$ws.Range("J65").Value = 4719
$ws.Range("K65").Value = 1267746.25
$ws.Range("L65").Value = 23595
$ws.Range("M65").Value = -1264626.25
$ws.Range("N65").Value = -29835

# Sheet ALC, row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2357.516
$ws.Range("J70").Value = 3026.5264
$ws.Range("L70").Value = 9079.5792
$ws.Range("N70").Value = -9619.5792

# Sheet ALC, row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2357.516
$ws.Range("J73").Value = 3026.5264
$ws.Range("L73").Value = 9079.5792
$ws.Range("N73").Value = -10951.5792

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2165.3333
$ws.Range("I86").Value = 2323
$ws.Range("K86").Value = 2323
$ws.Range("M86").Value = -1200

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2165.3333
$ws.Range("I89").Value = 2323
$ws.Range("K89").Value = 11615
$ws.Range("M89").Value = -5999

# Sheet ALC, row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 20043.104
$ws.Range("I92").Value = 1148.0541
$ws.Range("J92").Value = 83599.17999999999
$ws.Range("K92").Value = 1148.0541
$ws.Range("L92").Value = 83599.17999999999
$ws.Range("M92").Value = 99.94589999999994
$ws.Range("N92").Value = -86095.17999999999

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 145326
$ws.Range("I106").Value = 335260.66
$ws.Range("J106").Value = 2875
$ws.Range("K106").Value = 335260.66
$ws.Range("L106").Value = 2875
$ws.Range("M106").Value = -334629.66
$ws.Range("N106").Value = -4137

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2440.7437
$ws.Range("I132").Value = 2442.8108
$ws.Range("K132").Value = 7328.432400000001
$ws.Range("M132").Value = -4798.432400000001

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5077.2607
$ws.Range("I141").Value = 5262.591
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 15787.773
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = -10607.773
$ws.Range("N141").Value = -13360

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 17858720
$ws.Range("I45").Value = 31251018
$ws.Range("K45").Value = 31251018
$ws.Range("M45").Value = -31250641

# Sheet ARM, row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 69398
$ws.Range("J113").Value = 69398
$ws.Range("L113").Value = 69398
$ws.Range("N113").Value = -78076

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3180.3333
$ws.Range("J132").Value = 2803.3333
$ws.Range("L132").Value = 8409.999899999999
$ws.Range("N132").Value = -13469.9999

# Sheet CRP, row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147
$ws.Range("I7").Value = 147
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 147
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -34
$ws.Range("N7").ClearContents()

# Sheet CRP, row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 621.3333
$ws.Range("I10").Value = 621.3333
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 621.3333
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -482.3333
$ws.Range("N10").ClearContents()

# Sheet CRP, row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7813.8
$ws.Range("J62").Value = 7690.5713
$ws.Range("L62").Value = 7690.5713
$ws.Range("N62").Value = -8938.5713

# Sheet CRP, row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 7813.8
$ws.Range("J65").Value = 7690.5713
$ws.Range("L65").Value = 38452.85649999999
$ws.Range("N65").Value = -44692.85649999999

# Sheet CUL, row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 2021.909
$ws.Range("I8").Value = 2021.909
$ws.Range("K8").Value = 6065.727000000001
$ws.Range("M8").Value = -5926.727000000001

# Sheet CUL, row 93
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 14212.63
$ws.Range("J93").Value = 14212.63
$ws.Range("L93").Value = 42637.89
$ws.Range("N93").Value = -46381.89

# Sheet CUL, row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 195.64285
$ws.Range("I97").Value = 142.9
$ws.Range("K97").Value = 428.7
$ws.Range("M97").Value = 67.29999999999995

# Sheet CUL, row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 41667760
$ws.Range("I113").Value = 294.6
$ws.Range("K113").Value = 883.8000000000001
$ws.Range("M113").Value = 1286.2

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6119.5
$ws.Range("I70").Value = 6490
$ws.Range("K70").Value = 6490
$ws.Range("M70").Value = -6220

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6119.5
$ws.Range("I73").Value = 6490
$ws.Range("K73").Value = 6490
$ws.Range("M73").Value = -5554

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 39421.08
$ws.Range("I126").Value = 48393.6
$ws.Range("K126").Value = 145180.8
$ws.Range("M126").Value = -142710.8

# Sheet LTW, row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5749.25
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 5749.25
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5749.25
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -7247.25

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5749.25
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 5749.25
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 28746.25
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -36234.25

# Sheet LTW, row 109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 58121.43
$ws.Range("I109").Value = 58280
$ws.Range("J109").Value = 58095
$ws.Range("K109").Value = 58280
$ws.Range("L109").Value = 58095
$ws.Range("M109").Value = -56893
$ws.Range("N109").Value = -60869

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2724.3
$ws.Range("I132").Value = 2296.0312
$ws.Range("K132").Value = 6888.0936
$ws.Range("M132").Value = -4358.0936

# Sheet WVR, row 74
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 12081.25
$ws.Range("J74").Value = 10163
$ws.Range("L74").Value = 10163
$ws.Range("N74").Value = -12035

# Sheet WVR, row 77
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 12081.25
$ws.Range("J77").Value = 10163
$ws.Range("L77").Value = 30489
$ws.Range("N77").Value = -39849

# Sheet WVR, row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 21741522
$ws.Range("I126").Value = 27780084
$ws.Range("K126").Value = 83340252
$ws.Range("M126").Value = -83337782
